$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Both "ReLU" labels live inside the "Group 327" shape (top-level shape
# index 119 on this slide): GroupItems 3 ("Rectangle 330") and
# GroupItems 12 ("Rectangle 339"). Each needs its text changed from
# "ReLU" to "LReLU" plus a new trailing run " (-0.1)".
$group = $s.Shapes.Item(119)

$targets = @(3, 12)
foreach ($idx in $targets) {
    $shape = $group.GroupItems.Item($idx)
    $tr = $shape.TextFrame.TextRange

    # Rename "ReLU" -> "LReLU" in place (keeps the original run/formatting).
    $word = $tr.Characters(1, 4)
    $word.Text = "LReLU"

    # Append " (-0.1)" as its own run right after, picking up the
    # paragraph's trailing (endParaRPr) formatting.
    $suffix = $tr.Characters(6, 7)
    $suffix.Text = " (-0.1)"
    $suffix.Font.Size = 28
}
